$d = $word.ActiveDocument

$replacements = @(
    @("45÷6=7, 3", "21÷8=2, 5"),
    @("41÷3=13, 2", "63÷3=21, 0"),
    @("49÷5=9, 4", "73÷5=14, 3"),
    @("86÷9=9, 5", "71÷2=35, 1"),
    @("68÷7=9, 5", "18÷2=9, 0"),
    @("80÷4=20, 0", "32÷9=3, 5"),
    @("71÷8=8, 7", "53÷5=10, 3"),
    @("43÷4=10, 3", "74÷3=24, 2"),
    @("78÷8=9, 6", "70÷7=10, 0"),
    @("43÷8=5, 3", "35÷6=5, 5"),
    @("75÷5=15, 0", "76÷2=38, 0"),
    @("69÷9=7, 6", "25÷5=5, 0"),
    @("85÷3=28, 1", "79÷2=39, 1"),
    @("14÷9=1, 5", "14÷3=4, 2"),
    @("61÷3=20, 1", "83÷4=20, 3"),
    @("84÷8=10, 4", "55÷5=11, 0"),
    @("88÷6=14, 4", "46÷8=5, 6"),
    @("62÷6=10, 2", "87÷4=21, 3"),
    @("93÷4=23, 1", "61÷3=20, 1"),
    @("76÷8=9, 4", "23÷2=11, 1"),
    @("86÷3=28, 2", "56÷9=6, 2"),
    @("54÷5=10, 4", "60÷3=20, 0"),
    @("48÷6=8, 0", "29÷8=3, 5"),
    @("55÷4=13, 3", "37÷9=4, 1"),
    @("89÷2=44, 1", "21÷2=10, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
